$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value2 = 295.3913
$ws.Range("I53").Value2 = 451
$ws.Range("J53").Value2 = 152.75
$ws.Range("K53").Value2 = 451
$ws.Range("L53").Value2 = 152.75
$ws.Range("M53").Value2 = 186
$ws.Range("N53").Value2 = -1426.75

$ws.Range("H92").Value2 = 644.5
$ws.Range("I92").Value2 = 700.8570999999999
$ws.Range("K92").Value2 = 700.8570999999999
$ws.Range("M92").Value2 = 547.1429000000001

$ws.Range("H115").Value2 = 243
$ws.Range("I115").Value2 = 243
$ws.Range("K115").Value2 = 729
$ws.Range("M115").Value2 = 838

$ws.Range("H137").Value2 = 2964.0908
$ws.Range("I137").Value2 = 2223.3333
$ws.Range("J137").Value2 = 3853
$ws.Range("K137").Value2 = 6669.999899999999
$ws.Range("L137").Value2 = 11559
$ws.Range("M137").Value2 = -4119.999899999999
$ws.Range("N137").Value2 = -16659

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1588.75
$ws.Range("I2").Value2 = 1516.1666
$ws.Range("J2").Value2 = 1806.5
$ws.Range("K2").Value2 = 1516.1666
$ws.Range("L2").Value2 = 1806.5
$ws.Range("M2").Value2 = -1403.1666
$ws.Range("N2").Value2 = -2032.5

$ws.Range("H32").Value2 = 5697.0454
$ws.Range("I32").Value2 = 4666.744
$ws.Range("K32").Value2 = 4666.744
$ws.Range("M32").Value2 = -4379.744

$ws.Range("H116").Value2 = 1588.75
$ws.Range("I116").Value2 = 1516.1666
$ws.Range("J116").Value2 = 1806.5
$ws.Range("K116").Value2 = 1516.1666
$ws.Range("L116").Value2 = 1806.5
$ws.Range("M116").Value2 = 777.8334
$ws.Range("N116").Value2 = -6394.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1588.75
$ws.Range("I3").Value2 = 1516.1666
$ws.Range("J3").Value2 = 1806.5
$ws.Range("K3").Value2 = 1516.1666
$ws.Range("L3").Value2 = 1806.5
$ws.Range("M3").Value2 = -1402.1666
$ws.Range("N3").Value2 = -2034.5

$ws.Range("H105").Value2 = 2766.7144
$ws.Range("I105").Value2 = 2473.1538
$ws.Range("K105").Value2 = 2473.1538
$ws.Range("M105").Value2 = -726.1538

$ws.Range("H107").Value2 = 1259
$ws.Range("I107").Value2 = 567.55554
$ws.Range("J107").Value2 = 3333.3333
$ws.Range("K107").Value2 = 567.55554
$ws.Range("L107").Value2 = 3333.3333
$ws.Range("M107").Value2 = 1352.44446
$ws.Range("N107").Value2 = -7173.3333

$ws.Range("H134").Value2 = 6712.1
$ws.Range("I134").Value2 = 6712.1
$ws.Range("K134").Value2 = 20136.3
$ws.Range("M134").Value2 = -17601.3

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value2 = 696.6
$ws.Range("I2").Value2 = 836
$ws.Range("J2").Value2 = 636.8570999999999
$ws.Range("K2").Value2 = 836
$ws.Range("L2").Value2 = 636.8570999999999
$ws.Range("M2").Value2 = -723
$ws.Range("N2").Value2 = -862.8570999999999

$ws.Range("H5").Value2 = 202.22223
$ws.Range("I5").Value2 = 160.71428
$ws.Range("J5").Value2 = 347.5
$ws.Range("K5").Value2 = 160.71428
$ws.Range("L5").Value2 = 347.5
$ws.Range("M5").Value2 = -48.71428
$ws.Range("N5").Value2 = -571.5

$ws.Range("H11").Value2 = 629.5
$ws.Range("J11").Value2 = 275
$ws.Range("L11").Value2 = 275
$ws.Range("N11").Value2 = -555

$ws.Range("H12").Value2 = 5000
$ws.Range("I12").Value2 = 5000
$ws.Range("K12").Value2 = 5000
$ws.Range("M12").Value2 = -4830

$ws.Range("H16").Value2 = 7598.25
$ws.Range("I16").Value2 = 3596.75
$ws.Range("J16").Value2 = 11599.75
$ws.Range("K16").Value2 = 3596.75
$ws.Range("L16").Value2 = 11599.75
$ws.Range("M16").Value2 = -3309.75
$ws.Range("N16").Value2 = -12173.75

$ws.Range("H31").Value2 = 3653.3333
$ws.Range("I31").Value2 = 2522.5
$ws.Range("J31").Value2 = 5915
$ws.Range("K31").Value2 = 2522.5
$ws.Range("L31").Value2 = 5915
$ws.Range("M31").Value2 = -2227.5
$ws.Range("N31").Value2 = -6505

$ws.Range("H34").Value2 = 3653.3333
$ws.Range("I34").Value2 = 2522.5
$ws.Range("J34").Value2 = 5915
$ws.Range("K34").Value2 = 2522.5
$ws.Range("L34").Value2 = 5915
$ws.Range("M34").Value2 = -2320.5
$ws.Range("N34").Value2 = -6319

$ws.Range("H105").Value2 = 1034.2858
$ws.Range("I105").Value2 = 1021.61536
$ws.Range("J105").Value2 = 1199
$ws.Range("K105").Value2 = 1021.61536
$ws.Range("L105").Value2 = 1199
$ws.Range("M105").Value2 = 725.38464
$ws.Range("N105").Value2 = -4693

$ws.Range("H113").Value2 = 7598.25
$ws.Range("I113").Value2 = 3596.75
$ws.Range("J113").Value2 = 11599.75
$ws.Range("K113").Value2 = 3596.75
$ws.Range("L113").Value2 = 11599.75
$ws.Range("M113").Value2 = -1426.75
$ws.Range("N113").Value2 = -15939.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 505.7143
$ws.Range("I5").Value2 = 501.41177
$ws.Range("J5").Value2 = 524
$ws.Range("K5").Value2 = 1504.23531
$ws.Range("L5").Value2 = 1572
$ws.Range("M5").Value2 = -1392.23531
$ws.Range("N5").Value2 = -1796

$ws.Range("H131").Value2 = 2334.7058
$ws.Range("I131").Value2 = 1574.25
$ws.Range("J131").Value2 = 2568.6924
$ws.Range("K131").Value2 = 4722.75
$ws.Range("L131").Value2 = 7706.0772
$ws.Range("M131").Value2 = 317.25
$ws.Range("N131").Value2 = -17786.0772

$ws.Range("H135").Value2 = 505.7143
$ws.Range("I135").Value2 = 501.41177
$ws.Range("J135").Value2 = 524
$ws.Range("K135").Value2 = 4512.70593
$ws.Range("L135").Value2 = 4716
$ws.Range("M135").Value2 = -1977.70593
$ws.Range("N135").Value2 = -9786

$ws.Range("H137").Value2 = 10190
$ws.Range("J137").Value2 = 12237.5
$ws.Range("L137").Value2 = 36712.5
$ws.Range("N137").Value2 = -46912.5

$ws.Range("H139").Value2 = 2520
$ws.Range("I139").Value2 = 2520
$ws.Range("K139").Value2 = 7560
$ws.Range("M139").Value2 = -2420

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value2 = 20000
$ws.Range("I20").Value2 = 1000
$ws.Range("J20").Value2 = 39000
$ws.Range("K20").Value2 = 1000
$ws.Range("L20").Value2 = 39000
$ws.Range("M20").Value2 = -755
$ws.Range("N20").Value2 = -39490

$ws.Range("H24").Value2 = 0
$ws.Range("I24").Value2 = 0
$ws.Range("K24").Value2 = 0
$ws.Range("M24").ClearContents()

$ws.Range("H46").Value2 = 14888.777

$ws.Range("H57").Value2 = 21666.666
$ws.Range("I57").Value2 = 15000
$ws.Range("J57").Value2 = 25000
$ws.Range("K57").Value2 = 15000
$ws.Range("L57").Value2 = 25000
$ws.Range("M57").Value2 = -14180
$ws.Range("N57").Value2 = -26640

$ws.Range("H126").Value2 = 1250
$ws.Range("I126").Value2 = 1250
$ws.Range("J126").Value2 = 0
$ws.Range("K126").Value2 = 3750
$ws.Range("L126").Value2 = 0
$ws.Range("M126").Value2 = -1280
$ws.Range("N126").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value2 = 4653.4546
$ws.Range("I46").Value2 = 5000
$ws.Range("K46").Value2 = 5000
$ws.Range("M46").Value2 = -4812

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value2 = 5000
$ws.Range("J25").Value2 = 5000
$ws.Range("L25").Value2 = 5000
$ws.Range("N25").Value2 = -5586

$ws.Range("H31").Value2 = 20450
$ws.Range("I31").Value2 = 1000
$ws.Range("J31").Value2 = 39900
$ws.Range("K31").Value2 = 1000
$ws.Range("L31").Value2 = 39900
$ws.Range("M31").Value2 = -652
$ws.Range("N31").Value2 = -40596

$ws.Range("H92").Value2 = 38199.6
$ws.Range("J92").Value2 = 38199.6
$ws.Range("L92").Value2 = 38199.6
$ws.Range("N92").Value2 = -43191.6

$ws.Range("H132").Value2 = 750
$ws.Range("I132").Value2 = 750
$ws.Range("K132").Value2 = 2250
$ws.Range("M132").Value2 = 280
